# Regenerate the "K" column (column G) values in the save_data sheet.
# Per the commit message: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" — the recomputed K values replace
# the previous ones for every data row (rows 2-65), except rows 15 and 24
# whose recomputed value happened to match the original value already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newK = [ordered]@{
    2  = 1
    3  = 2
    4  = 3
    5  = 2
    6  = 2
    7  = 3
    8  = 2
    9  = 3
    10 = 4
    11 = 2
    12 = 2
    13 = 1
    14 = 2
    16 = 1
    17 = 1
    18 = 3
    19 = 2
    20 = 3
    21 = 2
    22 = 2
    23 = 0
    25 = 2
    26 = 1
    27 = 2
    28 = 2
    29 = 2
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 0
    40 = 2
    41 = 1
    42 = 1
    43 = 2
    44 = 0
    45 = 2
    46 = 2
    47 = 1
    48 = 1
    49 = 2
    50 = 2
    51 = 3
    52 = 2
    53 = 0
    54 = 2
    55 = 1
    56 = 1
    57 = 2
    58 = 1
    59 = 2
    60 = 2
    61 = 1
    62 = 1
    63 = 1
    64 = 0
    65 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
